$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 89
$ws.Range("I8").Value = 89
$ws.Range("K8").Value = 267
$ws.Range("M8").Value = -128
$ws.Range("H80").Value = 4344.222
$ws.Range("I80").Value = 5100
$ws.Range("K80").Value = 15300
$ws.Range("M80").Value = -14302
$ws.Range("H83").Value = 4344.222
$ws.Range("I83").Value = 5100
$ws.Range("K83").Value = 45900
$ws.Range("M83").Value = -40908
$ws.Range("H88").Value = 10000
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -9594
$ws.Range("H91").Value = 10000
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -8596
$ws.Range("H141").Value = 3095
$ws.Range("I141").Value = 3095
$ws.Range("K141").Value = 9285
$ws.Range("M141").Value = -4105

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 67
$ws.Range("I4").Value = 69.666664
$ws.Range("J4").Value = 59
$ws.Range("K4").Value = 69.666664
$ws.Range("L4").Value = 59
$ws.Range("M4").Value = 46.333336
$ws.Range("N4").Value = -291
$ws.Range("H45").Value = 2651.6667
$ws.Range("I45").Value = 2972.1
$ws.Range("K45").Value = 2972.1
$ws.Range("M45").Value = -2595.1
$ws.Range("H110").Value = 3566.8572
$ws.Range("I110").Value = 2135.6365
$ws.Range("J110").Value = 8814.666999999999
$ws.Range("K110").Value = 2135.6365
$ws.Range("L110").Value = 8814.666999999999
$ws.Range("M110").Value = -90.63650000000007
$ws.Range("N110").Value = -12904.667
$ws.Range("H132").Value = 2718.4614
$ws.Range("I132").Value = 1730.375
$ws.Range("K132").Value = 5191.125
$ws.Range("M132").Value = -2661.125

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 30000
$ws.Range("J10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("N10").Value = -30280
$ws.Range("H64").Value = 2221.7778
$ws.Range("I64").Value = 586.75
$ws.Range("K64").Value = 586.75
$ws.Range("M64").Value = -361.75
$ws.Range("H67").Value = 2221.7778
$ws.Range("I67").Value = 586.75
$ws.Range("K67").Value = 586.75
$ws.Range("M67").Value = 193.25
$ws.Range("H107").Value = 1373.5
$ws.Range("I107").Value = 1462.5714
$ws.Range("K107").Value = 1462.5714
$ws.Range("M107").Value = 457.4286
$ws.Range("H141").Value = 159995.5
$ws.Range("J141").Value = 159995
$ws.Range("L141").Value = 159995
$ws.Range("N141").Value = -170355

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1195.6666
$ws.Range("I16").Value = 1483
$ws.Range("J16").Value = 908.3333
$ws.Range("K16").Value = 1483
$ws.Range("L16").Value = 908.3333
$ws.Range("M16").Value = -1196
$ws.Range("N16").Value = -1482.3333
$ws.Range("H113").Value = 1195.6666
$ws.Range("I113").Value = 1483
$ws.Range("J113").Value = 908.3333
$ws.Range("K113").Value = 1483
$ws.Range("L113").Value = 908.3333
$ws.Range("M113").Value = 687
$ws.Range("N113").Value = -5248.3333
$ws.Range("H141").Value = 149992
$ws.Range("J141").Value = 149992
$ws.Range("L141").Value = 149992
$ws.Range("N141").Value = -160352

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1067.1
$ws.Range("I5").Value = 1459.3334
$ws.Range("J5").Value = 899
$ws.Range("K5").Value = 4378.0002
$ws.Range("L5").Value = 2697
$ws.Range("M5").Value = -4266.0002
$ws.Range("N5").Value = -2921
$ws.Range("H8").Value = 1334978.2
$ws.Range("I8").Value = 1334978.2
$ws.Range("K8").Value = 4004934.6
$ws.Range("M8").Value = -4004795.6
$ws.Range("H14").Value = 598.75
$ws.Range("I14").Value = 598.75
$ws.Range("K14").Value = 1796.25
$ws.Range("M14").Value = -1623.25
$ws.Range("H69").Value = 20000
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61622
$ws.Range("H72").Value = 20000
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -188112
$ws.Range("H75").Value = 4259.8
$ws.Range("J75").Value = 3074.75
$ws.Range("L75").Value = 9224.25
$ws.Range("N75").Value = -11220.25
$ws.Range("H78").Value = 4259.8
$ws.Range("J78").Value = 3074.75
$ws.Range("L78").Value = 27672.75
$ws.Range("N78").Value = -37656.75
$ws.Range("H92").Value = 508.16666
$ws.Range("J92").Value = 566.3333
$ws.Range("L92").Value = 1698.9999
$ws.Range("N92").Value = -4194.9999
$ws.Range("H103").Value = 397
$ws.Range("I103").Value = 295
$ws.Range("J103").Value = 499
$ws.Range("K103").Value = 885
$ws.Range("L103").Value = 1497
$ws.Range("M103").Value = -6
$ws.Range("N103").Value = -3255
$ws.Range("H113").Value = 1175.4615
$ws.Range("I113").Value = 320.5
$ws.Range("J113").Value = 1555.4445
$ws.Range("K113").Value = 961.5
$ws.Range("L113").Value = 4666.333500000001
$ws.Range("M113").Value = 1208.5
$ws.Range("N113").Value = -9006.333500000001
$ws.Range("H117").Value = 1011.6
$ws.Range("I117").Value = 329
$ws.Range("J117").Value = 1466.6666
$ws.Range("K117").Value = 987
$ws.Range("L117").Value = 4399.9998
$ws.Range("M117").Value = 2455
$ws.Range("N117").Value = -11283.9998
$ws.Range("H132").Value = 1332.8334
$ws.Range("I132").Value = 1332.8334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11995.5006
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -9465.500599999999
$ws.Range("H135").Value = 1067.1
$ws.Range("I135").Value = 1459.3334
$ws.Range("J135").Value = 899
$ws.Range("K135").Value = 13134.0006
$ws.Range("L135").Value = 8091
$ws.Range("M135").Value = -10599.0006
$ws.Range("N135").Value = -13161

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3780.3333
$ws.Range("I113").Value = 3780.3333
$ws.Range("K113").Value = 3780.3333
$ws.Range("M113").Value = -1610.3333
$ws.Range("H119").Value = 150000
$ws.Range("J119").Value = 150000
$ws.Range("L119").Value = 150000
$ws.Range("N119").Value = -159676

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2930.8333
$ws.Range("I7").Value = 2930.8333
$ws.Range("K7").Value = 2930.8333
$ws.Range("M7").Value = -2818.8333
$ws.Range("H26").Value = 50000
$ws.Range("J26").Value = 50000
$ws.Range("L26").Value = 50000
$ws.Range("N26").Value = -50590
$ws.Range("H82").Value = 3098.125
$ws.Range("I82").Value = 3200
$ws.Range("K82").Value = 3200
$ws.Range("M82").Value = -2839
$ws.Range("H85").Value = 3098.125
$ws.Range("I85").Value = 3200
$ws.Range("K85").Value = 3200
$ws.Range("M85").Value = -1952
$ws.Range("H126").Value = 2930.8333
$ws.Range("I126").Value = 2930.8333
$ws.Range("K126").Value = 8792.499899999999
$ws.Range("M126").Value = -6322.499899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 49990
$ws.Range("J74").Value = 49990
$ws.Range("L74").Value = 49990
$ws.Range("N74").Value = -51862
$ws.Range("H77").Value = 49990
$ws.Range("J77").Value = 49990
$ws.Range("L77").Value = 149970
$ws.Range("N77").Value = -159330
$ws.Range("H81").Value = 1868.875
$ws.Range("I81").Value = 1868.875
$ws.Range("K81").Value = 3737.75
$ws.Range("M81").Value = -2676.75
$ws.Range("H84").Value = 1868.875
$ws.Range("I84").Value = 1868.875
$ws.Range("K84").Value = 18688.75
$ws.Range("M84").Value = -13384.75
$ws.Range("H100").Value = 3224.8333
$ws.Range("I100").Value = 2924.5
$ws.Range("J100").Value = 3375
$ws.Range("K100").Value = 5849
$ws.Range("L100").Value = 6750
$ws.Range("M100").Value = -5308
$ws.Range("N100").Value = -7832
